# Applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.479.22"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.573.37"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3664"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.134"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07454"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.009"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.935"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "1.578.73"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001112"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06754"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.419"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "22.482.12"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.387"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.627"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").Value = "1.753.89"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.047"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.197"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.000"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.883"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08295"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02448"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2270"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06478"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.468"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.300"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.78%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6352"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6117"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.764"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.059"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.225"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07241"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.64%  "
